$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")

# Add the hyperlink + new GitHub-link text to F4 (the new "Square root of N
# upto 3 decimal places" row's Github Link cell).
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/2_square_root_of_n_upto_3_decimal_places.java", "", "", "dsa/2_square_root_of_n_upto_3_decimal_places.java at main " + [char]0xB7 + " ankurnecessary/dsa " + [char]0xB7 + " GitHub")

# Hyperlinks.Add() re-derives the cell style from scratch; bring it back in
# line with the existing hyperlink-styled cell (F3) so no new style record
# is introduced.
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)

# Row 4 grows taller to fit the wrapped text in the new F4 cell.
$ws.Rows("4").RowHeight = 57.6

# Selection moves on from F4 to B5 after the edit.
$ws.Range("B5").Select() | Out-Null
